$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lesson 3 (rows 27-36) fill-in-the-blank exercise data.
# Each row: B = Ukrainian sentence, D = English translation.
$rows = @(
    @{ Row = 27; B = "Він завжди миє свої руки перед обідом"; D = "He always washes his hands before dinner" },
    @{ Row = 28; B = "Вони допомагають своїм друзям"; D = "They help their friends" },
    @{ Row = 29; B = "Я ходжу в зал по понеділкам"; D = "I go to the gym on Mondays" },
    @{ Row = 30; B = "Вона працює в магазині"; D = "She works in the store" },
    @{ Row = 31; B = "Її син любить гратися"; D = "Her son loves to play" },
    @{ Row = 32; B = "Мій батько бігає вранці"; D = "My father runs in the mornings" },
    @{ Row = 33; B = "У йьому місці подають гарну каву"; D = "This place serves good coffee" },
    @{ Row = 34; B = "Мені потрібна інша ручка"; D = "I need another pen" },
    @{ Row = 35; B = "Ви мусите подумати про здоров'я"; D = "You must think about your health" },
    @{ Row = 36; B = "Ми їздимо на тому ж автобусі кожного ранку"; D = "We take the same bus every morning" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}

# Update the view/selection state left by the editing session.
$ws.Range("H30").Select()
